# Update the "cryptos" listing with refreshed prices/volume figures
# (and two rows that swapped rank position), per the scheduled
# GitHub Actions refresh.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rows whose rank (and therefore Coin/Link) stayed the same, only
#     Price (D) and/or Volume(1h) (E) changed. $null means "leave D as-is".
$updates = @(
    @{ Row = 2;  D = "94.893.48";   E = "  -1.63%  " },
    @{ Row = 3;  D = "3.563.99";    E = "  -1.64%  " },
    @{ Row = 4;  D = $null;         E = "  -0.03%  " },
    @{ Row = 5;  D = "236.02";      E = "  -2.38%  " },
    @{ Row = 6;  D = "654.11";      E = "  +2.07%  " },
    @{ Row = 7;  D = $null;         E = "  -1.11%  " },
    @{ Row = 8;  D = $null;         E = "  -1.38%  " },
    @{ Row = 9;  D = $null;         E = "  +0.11%  " },
    @{ Row = 10; D = $null;         E = "  -1.24%  " },
    @{ Row = 11; D = "3.563.25";    E = "  -1.54%  " },
    @{ Row = 12; D = $null;         E = "  +0.60%  " },
    @{ Row = 13; D = "42.20";       E = "  -2.52%  " },
    @{ Row = 14; D = "6.44";        E = "  +1.06%  " },
    @{ Row = 15; D = "4.224.02";    E = "  -1.80%  " },
    @{ Row = 16; D = "94.738.15";   E = "  -1.71%  " },
    @{ Row = 17; D = $null;         E = "  -0.88%  " },
    @{ Row = 20; D = $null;         E = "  -4.13%  " },
    @{ Row = 21; D = "17.71";       E = "  -3.08%  " },
    @{ Row = 22; D = $null;         E = "  -0.63%  " },
    @{ Row = 23; D = "507.17";      E = "  -1.77%  " },
    @{ Row = 24; D = "0.480";       E = "  -4.47%  " },
    @{ Row = 25; D = $null;         E = "  +0.53%  " },
    @{ Row = 26; D = "0.0000195";   E = "  -2.14%  " },
    @{ Row = 27; D = "94.80";       E = "  -3.52%  " },
    @{ Row = 28; D = "12.54";       E = "  +0.00%  " },
    @{ Row = 29; D = "3.753.81";    E = "  -1.48%  " },
    @{ Row = 30; D = $null;         E = "  -5.49%  " },
    @{ Row = 31; D = $null;         E = "  -0.68%  " },
    @{ Row = 32; D = "11.45";       E = "  -1.60%  " },
    @{ Row = 33; D = $null;         E = "  +0.10%  " },
    @{ Row = 34; D = "1.00";        E = "  +0.50%  " },
    @{ Row = 35; D = $null;         E = "  -3.87%  " },
    @{ Row = 36; D = "31.72";       E = "  +3.76%  " },
    @{ Row = 37; D = $null;         E = "  +14.61%  " },
    @{ Row = 38; D = "0.554";       E = "  -2.76%  " },
    @{ Row = 39; D = "8.45";        E = "  +7.22%  " },
    @{ Row = 40; D = "578.75";      E = "  +0.20%  " },
    @{ Row = 41; D = $null;         E = "  +0.09%  " },
    @{ Row = 42; D = $null;         E = "  -1.20%  " },
    @{ Row = 43; D = "0.903";       E = "  -2.38%  " },
    @{ Row = 44; D = "1.81";        E = "  +3.78%  " },
    @{ Row = 47; D = $null;         E = "  +2.35%  " },
    @{ Row = 48; D = $null;         E = "  -1.77%  " },
    @{ Row = 49; D = $null;         E = "  -4.75%  " },
    @{ Row = 50; D = $null;         E = "  +0.48%  " },
    @{ Row = 51; D = "8.13";        E = "  -1.07%  " }
)

# Price column cells whose new text would otherwise be auto-parsed by Excel
# as a plain number (single decimal point, no thousands separators). Force
# them to stay plain text -- same as the source data -- by pre-formatting
# as Text before assigning the value, then restoring the original "Normal"
# style so no visible formatting change is left behind. (Rows 18/45 are
# handled in the rank-swap section below and use the same treatment.)
$textForceRows = @(5,6,13,14,21,23,24,26,27,28,32,34,36,38,39,40,43,44,51)

foreach ($u in $updates) {
    if ($null -ne $u.D) {
        $cell = $ws.Cells.Item($u.Row, 4)
        if ($textForceRows -contains $u.Row) {
            $cell.NumberFormat = "@"
            $cell.Value = $u.D
            $cell.Style = "Normal"
        } else {
            $cell.Value = $u.D
        }
    }
    $ws.Cells.Item($u.Row, 5).Value = $u.E
}

# --- Rows 18/19: WrappedEther and Polkadot swap rank order (Polkadot
#     jumped ahead of WrappedEther).
$ws.Cells.Item(18, 2).Value = "Polkadot"
$ws.Cells.Item(18, 3).Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$d18 = $ws.Cells.Item(18, 4)
$d18.NumberFormat = "@"
$d18.Value = "8.49"
$d18.Style = "Normal"
$ws.Cells.Item(18, 5).Value = "  +5.67%  "

$ws.Cells.Item(19, 2).Value = "WrappedEther"
$ws.Cells.Item(19, 3).Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Cells.Item(19, 4).Value = "3.565.58"
$ws.Cells.Item(19, 5).Value = "  -1.61%  "

# --- Rows 45/46: EnergySwap and Filecoin swap rank order (EnergySwap
#     jumped ahead of Filecoin).
$ws.Cells.Item(45, 2).Value = "EnergySwap"
$ws.Cells.Item(45, 3).Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$d45 = $ws.Cells.Item(45, 4)
$d45.NumberFormat = "@"
$d45.Value = "34.60"
$d45.Style = "Normal"
$ws.Cells.Item(45, 5).Value = "  +31.39%  "

$ws.Cells.Item(46, 2).Value = "Filecoin"
$ws.Cells.Item(46, 3).Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$d46 = $ws.Cells.Item(46, 4)
$d46.NumberFormat = "@"
$d46.Value = "5.73"
$d46.Style = "Normal"
$ws.Cells.Item(46, 5).Value = "  +0.22%  "
